$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated filenames (TC01 -> TC13, DNBSEQ-G400 -> NotSpecifiedindata)
$neo4jFile = @'
TC13_CDS_Filter_InstrumentModel-NotSpecifiedindata_Neo4jData.xlsx
'@
$webFile = @'
TC13_CDS_Filter_InstrumentModel-NotSpecifiedindata_WebData.xlsx
'@

# Updated queries (DNBSEQ-G400 -> Not specified in data)
$statQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@
$participantsQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$samplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$filesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Not specified in data']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

# Row 2 = ParticipantsTab
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# Row 3 = SamplesTab
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# Row 4 = FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# Update column widths for D and E (values tuned so the engine's internal
# character-width quantization lands on the target serialized width)
$ws.Columns.Item(4).ColumnWidth = 93.66
$ws.Columns.Item(5).ColumnWidth = 92.17

# Update the selected cell in the sheet view
$ws.Range("D2").Select()
